$d = $word.ActiveDocument

# --- Change 1 & 2: "Started." -> "Yes." in the two "Dealt With" cells ---
# (Find.Execute's Replace mode in this runtime operates over the whole
#  document story rather than the scoping Range, so we set the cell's
#  Range.Text directly instead — that is properly scoped to the cell.)
$d.Tables.Item(1).Cell(5, 3).Range.Text = "Yes."
$d.Tables.Item(2).Cell(6, 3).Range.Text = "Yes."

# --- Change 3: split the "vaccine passport" story's run, wrapping "so" in
#     w:proofErr gramStart/gramEnd markers (mirroring a Word grammar-check
#     pass). The Word object model has no direct way to author
#     w:proofErr, so we rebuild the paragraph via Range.InsertXML with the
#     exact paragraph/run identity (paraId/textId/rsids) preserved. ---
$cell = $d.Tables.Item(8).Cell(2, 2)
$start = $cell.Range.Start
$end = $cell.Range.End - 1
$rng = $d.Range($start, $end)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0DC79678" w14:textId="5B17AB01" w:rsidR="00E85D95" w:rsidRDefault="00E85D95" w:rsidP="002971E3"><w:r w:rsidRPr="00454F70"><w:t xml:space="preserve">As a player, I don&#8217;t want the actual code recorded. E.g., </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00454F70"><w:t>so</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00454F70"><w:t xml:space="preserve"> I can scan and score my vaccine passport.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)
